# Auto-generated edit script: normalizes "开始-结束" time range separators to "开始 - 结束"
# and updates 想去人数 (F) / 最低票价 (G) counters per the target diff.
$wb = $excel.ActiveWorkbook

# --- Worksheet 1 (sheet1) ---
$ws = $wb.Worksheets.Item(1)
$updates = @(
    @{ Row=2; E='2024.03.30 10:00 - 04.29 22:00'; F=253 }
    @{ Row=3; E='2024.04.12 10:00 - 05.12 20:00'; F=865 }
    @{ Row=4; E='2024.04.13 11:00 - 04.14 17:00' }
    @{ Row=5; E='2024.04.13 10:30 - 04.14 16:30'; F=2293 }
    @{ Row=6; E='2024.04.13 11:00 - 04.14 18:00'; F=1366 }
    @{ Row=7; E='2024.04.13 10:00 - 04.14 18:00'; F=119 }
    @{ Row=8; E='2024.04.13 10:00 - 04.13 17:00'; F=814 }
    @{ Row=9; E='2024.04.13 10:00 - 04.21 17:00'; F=1160 }
    @{ Row=10; E='2024.04.20 10:00 - 04.21 17:00'; F=1045 }
    @{ Row=11; E='2024.04.20 10:00 - 04.21 17:00'; F=3036 }
    @{ Row=12; E='2024.04.20 10:00 - 04.21 17:00' }
    @{ Row=13; E='2024.04.20 13:50 - 04.20 18:00'; F=46 }
    @{ Row=14; E='2024.04.20 10:00 - 04.21 17:00'; F=1105 }
    @{ Row=15; E='2024.04.20 10:00 - 04.21 17:00'; F=618 }
    @{ Row=16; E='2024.04.27 12:40 - 04.27 16:40' }
    @{ Row=17; E='2024.04.27 10:00 - 04.27 18:00' }
    @{ Row=18; E='2024.04.27 10:30 - 04.27 16:30' }
    @{ Row=19; E='2024.05.01 10:00 - 06.02 22:00'; F=1105 }
    @{ Row=20; E='2024.05.01 10:00 - 06.02 22:00'; F=1105 }
    @{ Row=21; E='2024.05.01 10:00 - 05.02 17:00'; F=155 }
    @{ Row=22; E='2024.05.01 10:30 - 05.01 16:30'; F=536 }
    @{ Row=23; E='2024.05.01 10:00 - 05.05 16:00'; F=184 }
    @{ Row=24; E='2024.05.01 10:00 - 05.01 17:00' }
    @{ Row=25; E='2024.05.02 10:30 - 05.04 19:00'; F=239 }
    @{ Row=26; E='2024.05.02 13:00 - 05.02 18:00'; F=645 }
    @{ Row=27; E='2024.05.02 10:20 - 05.03 16:30'; F=606 }
    @{ Row=28; E='2024.05.03 10:00 - 05.04 16:00'; F=5 }
    @{ Row=29; E='2024.05.04 10:00 - 05.05 17:00'; F=839 }
    @{ Row=30; E='2024.05.04 10:00 - 05.04 18:00'; F=75 }
    @{ Row=31; E='2024.05.05 10:00 - 05.05 17:00' }
    @{ Row=32; E='2024.05.05 10:00 - 05.05 18:00'; F=43 }
    @{ Row=33; E='2024.05.18 10:00 - 05.19 17:00'; F=1045 }
    @{ Row=34; E='2024.05.18 10:00 - 05.19 17:00'; F=5051 }
    @{ Row=35; E='2024.05.18 10:00 - 05.18 17:00'; F=518 }
    @{ Row=36; E='2024.05.19 10:00 - 05.19 17:00'; F=249 }
    @{ Row=37; E='2024.06.01 10:00 - 06.02 17:00'; F=136 }
    @{ Row=38; E='2024.06.08 10:00 - 06.10 16:00'; F=4 }
    @{ Row=39; E='2024.07.05 10:00 - 07.07 16:00' }
    @{ Row=40; E='2024.07.12 10:00 - 07.14 16:00'; F=3 }
)
foreach ($u in $updates) {
    if ($u.ContainsKey("E")) { $ws.Cells.Item($u.Row, 5).Value = $u.E }
    if ($u.ContainsKey("F")) { $ws.Cells.Item($u.Row, 6).Value = $u.F }
    if ($u.ContainsKey("G")) { $ws.Cells.Item($u.Row, 7).Value = $u.G }
}

# --- Worksheet 2 (sheet2) ---
$ws = $wb.Worksheets.Item(2)
$updates = @(
    @{ Row=2; E='2024.04.05 15:50 - 05.03 20:10' }
    @{ Row=3; E='2024.04.12 19:30 - 04.12 21:30'; F=23 }
    @{ Row=4; E='2024.04.12 19:00 - 04.12 20:30'; F=364 }
    @{ Row=5; E='2024.04.13 16:00 - 04.13 18:00' }
    @{ Row=6; E='2024.04.13 19:30 - 04.13 21:30'; F=410 }
    @{ Row=7; E='2024.04.13 19:30 - 04.13 21:00' }
    @{ Row=8; E='2024.04.13 18:30 - 04.20 21:50' }
    @{ Row=9; E='2024.04.14 14:00 - 04.14 15:40' }
    @{ Row=10; E='2024.04.14 15:00 - 05.01 20:15' }
    @{ Row=11; E='2024.04.20 19:30 - 04.20 21:00'; F=196; G=580 }
    @{ Row=12; E='2024.04.20 19:30 - 04.20 21:30' }
    @{ Row=13; E='2024.04.20 13:00 - 04.20 15:00'; F=289 }
    @{ Row=14; E='2024.04.21 15:00 - 04.21 16:30' }
    @{ Row=15; E='2024.04.21 20:00 - 04.21 21:30' }
    @{ Row=16; E='2024.04.25 19:30 - 04.25 21:00' }
    @{ Row=17; E='2024.04.26 19:30 - 04.26 21:30' }
    @{ Row=18; E='2024.04.26 19:00 - 04.26 20:30' }
    @{ Row=19; E='2024.04.27 20:00 - 04.27 21:30' }
    @{ Row=20; E='2024.05.01 20:00 - 05.01 22:00' }
    @{ Row=21; E='2024.05.01 18:30 - 05.01 21:00'; F=42 }
    @{ Row=22; E='2024.05.01 19:30 - 05.19 21:00'; F=302 }
    @{ Row=23; E='2024.05.02 13:30 - 05.02 15:20' }
    @{ Row=24; E='2024.05.02 19:30 - 05.02 21:00'; F=50 }
    @{ Row=25; E='2024.05.03 19:00 - 05.03 22:00'; F=381 }
    @{ Row=26; E='2024.05.03 19:30 - 05.03 21:00' }
    @{ Row=27; E='2024.05.04 14:20 - 06.09 15:35' }
    @{ Row=28; E='2024.05.04 14:00 - 05.04 16:00'; F=681 }
    @{ Row=29; E='2024.05.04 19:30 - 05.04 21:00' }
    @{ Row=30; E='2024.05.04 16:30 - 06.02 17:50' }
    @{ Row=31; E='2024.05.05 13:00 - 05.05 15:30' }
    @{ Row=32; E='2024.05.17 19:30 - 05.17 21:00' }
    @{ Row=33; E='2024.05.18 19:30 - 05.18 21:00' }
    @{ Row=34; E='2024.05.18 14:00 - 05.18 20:30' }
    @{ Row=35; E='2024.05.19 19:30 - 05.19 21:00' }
    @{ Row=36; E='2024.05.19 14:30 - 05.19 16:00' }
    @{ Row=37; E='2024.05.19 14:00 - 05.19 15:30'; F=441 }
    @{ Row=38; E='2024.05.24 19:30 - 05.24 21:00' }
    @{ Row=39; E='2024.05.25 19:30 - 05.25 21:00'; F=11 }
    @{ Row=40; E='2024.06.01 19:30 - 06.01 21:00' }
    @{ Row=41; E='2024.06.01 19:30 - 06.01 21:00' }
    @{ Row=42; E='2024.06.08 19:30 - 06.08 21:00' }
    @{ Row=43; E='2024.06.15 19:30 - 06.15 22:00' }
    @{ Row=44; E='2024.06.22 19:30 - 06.22 21:30' }
    @{ Row=45; E='2024.07.17 19:30 - 07.17 21:00'; F=3 }
    @{ Row=46; E='2024.07.19 19:30 - 07.19 21:30' }
)
foreach ($u in $updates) {
    if ($u.ContainsKey("E")) { $ws.Cells.Item($u.Row, 5).Value = $u.E }
    if ($u.ContainsKey("F")) { $ws.Cells.Item($u.Row, 6).Value = $u.F }
    if ($u.ContainsKey("G")) { $ws.Cells.Item($u.Row, 7).Value = $u.G }
}

# --- Worksheet 3 (sheet3) ---
$ws = $wb.Worksheets.Item(3)
$updates = @(
    @{ Row=2; E='2023.10.16 10:00 - 2024.10.15 21:00' }
    @{ Row=3; E='2023.10.25 10:00 - 2024.10.20 21:00' }
    @{ Row=4; E='2024.03.21 00:00 - 04.28 23:59'; F=642 }
    @{ Row=5; E='2024.04.04 00:00 - 05.20 23:59'; F=432 }
    @{ Row=6; E='2024.04.24 00:00 - 05.22 23:59'; F=415 }
)
foreach ($u in $updates) {
    if ($u.ContainsKey("E")) { $ws.Cells.Item($u.Row, 5).Value = $u.E }
    if ($u.ContainsKey("F")) { $ws.Cells.Item($u.Row, 6).Value = $u.F }
    if ($u.ContainsKey("G")) { $ws.Cells.Item($u.Row, 7).Value = $u.G }
}

# --- Worksheet 4 (sheet4) ---
$ws = $wb.Worksheets.Item(4)
$updates = @(
    @{ Row=2; E='2024.03.21 00:00 - 04.28 23:59'; F=642 }
    @{ Row=3; E='2024.03.30 10:00 - 04.29 22:00'; F=253 }
    @{ Row=4; E='2024.04.04 00:00 - 05.20 23:59'; F=432 }
    @{ Row=5; E='2024.04.12 10:00 - 05.12 20:00'; F=865 }
    @{ Row=6; E='2024.04.13 19:30 - 04.13 21:30'; F=410 }
    @{ Row=7; E='2024.04.13 11:00 - 04.14 17:00' }
    @{ Row=8; E='2024.04.13 10:30 - 04.14 16:30'; F=2293 }
    @{ Row=9; E='2024.04.13 11:00 - 04.14 18:00'; F=1366 }
    @{ Row=10; E='2024.04.13 10:00 - 04.14 18:00'; F=119 }
    @{ Row=11; E='2024.04.13 10:00 - 04.13 17:00'; F=814 }
    @{ Row=12; E='2024.04.13 10:00 - 04.21 17:00'; F=1160 }
    @{ Row=13; E='2024.04.14 14:00 - 04.14 15:40' }
    @{ Row=14; E='2024.04.20 19:30 - 04.20 21:00'; F=196; G=580 }
    @{ Row=15; E='2024.04.20 10:00 - 04.21 17:00'; F=1045 }
    @{ Row=16; E='2024.04.20 10:00 - 04.21 17:00'; F=3036 }
    @{ Row=17; E='2024.04.20 10:00 - 04.21 17:00' }
    @{ Row=18; E='2024.04.20 13:00 - 04.20 15:00'; F=289 }
    @{ Row=19; E='2024.04.20 10:00 - 04.21 17:00'; F=1105 }
    @{ Row=20; E='2024.04.20 10:00 - 04.21 17:00'; F=618 }
    @{ Row=21; E='2024.04.24 00:00 - 05.22 23:59'; F=415 }
    @{ Row=22; E='2024.04.27 12:40 - 04.27 16:40' }
    @{ Row=23; E='2024.04.27 10:00 - 04.27 18:00' }
    @{ Row=24; E='2024.04.27 10:30 - 04.27 16:30' }
    @{ Row=25; E='2024.05.01 10:00 - 06.02 22:00'; F=1105 }
    @{ Row=26; E='2024.05.01 10:00 - 06.02 22:00'; F=1105 }
    @{ Row=27; E='2024.05.01 10:00 - 05.02 17:00'; F=155 }
    @{ Row=28; E='2024.05.01 18:30 - 05.01 21:00'; F=42 }
    @{ Row=29; E='2024.05.01 10:30 - 05.01 16:30'; F=536 }
    @{ Row=30; E='2024.05.01 19:30 - 05.19 21:00'; F=302 }
    @{ Row=31; E='2024.05.01 10:00 - 05.05 16:00'; F=184 }
    @{ Row=32; E='2024.05.02 10:30 - 05.04 19:00'; F=239 }
    @{ Row=33; E='2024.05.02 19:30 - 05.02 21:00'; F=50 }
    @{ Row=34; E='2024.05.02 13:00 - 05.02 18:00'; F=645 }
    @{ Row=35; E='2024.05.02 10:20 - 05.03 16:30'; F=606 }
    @{ Row=36; E='2024.05.03 19:00 - 05.03 22:00'; F=381 }
    @{ Row=37; E='2024.05.04 14:00 - 05.04 16:00'; F=681 }
    @{ Row=38; E='2024.05.04 10:00 - 05.05 17:00' }
    @{ Row=39; E='2024.05.04 19:30 - 05.04 21:00' }
    @{ Row=40; E='2024.05.04 10:00 - 05.04 18:00'; F=75 }
    @{ Row=41; E='2024.05.05 10:00 - 05.05 17:00' }
    @{ Row=42; E='2024.05.05 13:00 - 05.05 15:30' }
    @{ Row=43; E='2024.05.18 10:00 - 05.19 17:00'; F=1045 }
    @{ Row=44; E='2024.05.18 10:00 - 05.19 17:00'; F=5051 }
    @{ Row=45; E='2024.05.18 14:00 - 05.18 20:30' }
    @{ Row=46; E='2024.05.18 10:00 - 05.18 17:00'; F=518 }
    @{ Row=47; E='2024.05.19 14:00 - 05.19 15:30'; F=441 }
    @{ Row=48; E='2024.05.19 14:00 - 05.19 15:30'; F=441 }
    @{ Row=49; E='2024.05.19 10:00 - 05.19 17:00'; F=249 }
    @{ Row=50; E='2024.06.08 10:00 - 06.10 16:00'; F=4 }
    @{ Row=51; E='2024.06.22 19:30 - 06.22 21:30' }
)
foreach ($u in $updates) {
    if ($u.ContainsKey("E")) { $ws.Cells.Item($u.Row, 5).Value = $u.E }
    if ($u.ContainsKey("F")) { $ws.Cells.Item($u.Row, 6).Value = $u.F }
    if ($u.ContainsKey("G")) { $ws.Cells.Item($u.Row, 7).Value = $u.G }
}

